$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.309509754180908
$ws.Range("B1").Value = 1.512053728103638
$ws.Range("C1").Value = 1.277461051940918
$ws.Range("D1").Value = 1.329763054847717
$ws.Range("E1").Value = 1.042883515357971
